# Update cryptos list values/links per latest data pull (Sat Nov 23 14:43:55 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.712.15"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "3.490.32"
$ws.Range("E3").Value = "  +5.70%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "262.15"
$ws.Range("E5").Value = "  +2.57%  "

$ws.Range("D6").Value = "678.82"
$ws.Range("E6").Value = "  +9.34%  "

$ws.Range("E7").Value = "  +9.32%  "

$ws.Range("E8").Value = "  +16.66%  "

$ws.Range("E9").Value = "  +22.53%  "

$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("D11").Value = "3.493.66"

$ws.Range("D12").Value = "0.222"
$ws.Range("E12").Value = "  +11.99%  "

$ws.Range("D13").Value = "42.71"
$ws.Range("E13").Value = "  +10.59%  "

$ws.Range("E14").Value = "  +11.54%  "

$ws.Range("D15").Value = "6.18"
$ws.Range("E15").Value = "  +13.48%  "

$ws.Range("D16").Value = "98.427.61"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "4.142.01"
$ws.Range("E17").Value = "  +5.66%  "

$ws.Range("E18").Value = "  +33.27%  "

$ws.Range("D19").Value = "3.484.17"
$ws.Range("E19").Value = "  +5.47%  "

$ws.Range("D20").Value = "17.62"
$ws.Range("E20").Value = "  +16.67%  "

$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").Value = "537.73"
$ws.Range("E22").Value = "  +12.59%  "

$ws.Range("D23").Value = "10.81"
$ws.Range("E23").Value = "  +14.89%  "

$ws.Range("E24").Value = "  +8.99%  "

$ws.Range("D25").Value = "0.451"
$ws.Range("E25").Value = "  +53.74%  "

$ws.Range("D26").Value = "6.41"
$ws.Range("E26").Value = "  +14.95%  "

$ws.Range("D27").Value = "102.65"
$ws.Range("E27").Value = "  +17.06%  "

$ws.Range("D28").Value = "13.01"
$ws.Range("E28").Value = "  +9.94%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "0.151"
$ws.Range("E29").Value = "  +15.66%  "

$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").Value = "0.201"
$ws.Range("E30").Value = "  +7.95%  "

$ws.Range("D31").Value = "11.41"
$ws.Range("E31").Value = "  +16.17%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("B33").Value = "PolygonEcosystemToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D33").Value = "0.588"
$ws.Range("E33").Value = "  +29.67%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "30.93"
$ws.Range("E34").Value = "  +12.52%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.960"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "2.24"
$ws.Range("E36").Value = "  +16.03%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "8.02"
$ws.Range("E37").Value = "  +12.07%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.161"
$ws.Range("E38").Value = "  +9.65%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "541.57"
$ws.Range("E39").Value = "  +10.63%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.44"
$ws.Range("E40").Value = "  +16.33%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "24.77"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.873"
$ws.Range("E42").Value = "  +9.70%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0439"
$ws.Range("E43").Value = "  +35.64%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.53"
$ws.Range("E44").Value = "  +12.40%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "3.77"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "8.26"
$ws.Range("E46").Value = "  +16.86%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").Value = "2.13"
$ws.Range("E48").Value = "  +12.80%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  +18.78%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "5.30"
$ws.Range("E50").Value = "  +14.58%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "51.41"
$ws.Range("E51").Value = "  +13.31%  "
